# Weekly update: insert two new price-report rows (Primera / Segunda) for
# Betarraga at Terminal Hortofrutícola Agro Chillán, dated 2023-05-31 (serial 45077),
# right after the existing row for 2022-01-26 (row 559). This pushes all
# subsequent rows down by two, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 560-561 (existing data from row 560 onward shifts down by 2)
$ws.Range("A560:R561").Insert()

# Row 560: Primera
$ws.Range("A560").Value = 7
$ws.Range("B560").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C560").Value = "Ñuble"
$ws.Range("D560").Value = 45077
$ws.Range("E560").Value = 16
$ws.Range("F560").Value = 100114014
$ws.Range("G560").Value = "Betarraga"
$ws.Range("H560").Value = "Sin especificar"
$ws.Range("I560").Value = "Primera"
$ws.Range("J560").Value = 200
$ws.Range("K560").Value = 900
$ws.Range("L560").Value = 900
$ws.Range("M560").Value = 900
$ws.Range("N560").Value = "$/paquete 5 unidades"
$ws.Range("O560").Value = "Provincia de Diguillín"
$ws.Range("P560").Value = 180
$ws.Range("Q560").Value = 5
$ws.Range("R560").Value = "Hortaliza"

# Row 561: Segunda
$ws.Range("A561").Value = 7
$ws.Range("B561").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C561").Value = "Ñuble"
$ws.Range("D561").Value = 45077
$ws.Range("E561").Value = 16
$ws.Range("F561").Value = 100114014
$ws.Range("G561").Value = "Betarraga"
$ws.Range("H561").Value = "Sin especificar"
$ws.Range("I561").Value = "Segunda"
$ws.Range("J561").Value = 200
$ws.Range("K561").Value = 700
$ws.Range("L561").Value = 700
$ws.Range("M561").Value = 700
$ws.Range("N561").Value = "$/paquete 5 unidades"
$ws.Range("O561").Value = "Provincia de Diguillín"
$ws.Range("P561").Value = 140
$ws.Range("Q561").Value = 5
$ws.Range("R561").Value = "Hortaliza"
